# ---------------------------------------------------------------------------
# Update "Co so du lieu do an 2.xlsx" DB-schema sheet:
#   - classes table: reorder status/teacher_id/subject_id tail
#   - add a new "schedule" table (id, class_id, period, date, time, is_hoc_bu)
#     laid out in column C alongside the classes table
#   - attendances table: add period + date fields, re-style two rows
#   - attendance_class table: shifts up one row, field status -> is_present
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- capture the "last row / white fill" style (currently only on A24)
# before we repurpose that cell, so we can re-apply the same style at its
# new home (A23).
$ws.Range("A24").Copy()
$ws.Range("A23").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ===========================================================================
# classes table (column A) - rows 16-23
# ===========================================================================
$ws.Range("A21").Value = "status (hoàn thành|chưa đủ số lượng|đang tiến trình|...)"
$ws.Range("A22").Value = "teacher_id"
$ws.Range("A23").Value = "subject_id"

# row 24 in column A no longer holds data - clear value + formatting
$ws.Range("A24").Clear()

# ===========================================================================
# schedule table (column C) - rows 22-28 (new table, placed right under
# the subscriptions table which already occupies C16:C20)
# ===========================================================================
$ws.Range("C22").Value = "schedule"
$ws.Range("C16").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C23").Value = "id"
$ws.Range("C17").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C24").Value = "class_id"
$ws.Range("A18").Copy()
$ws.Range("C24").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C25").Value = "period"
$ws.Range("A18").Copy()
$ws.Range("C25").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C26").Value = "date"
$ws.Range("A18").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C27").Value = "time"
$ws.Range("A18").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C28").Value = "is_hoc_bu"
$ws.Range("A9").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ===========================================================================
# attendances table (column A) - rows 26-30
# ===========================================================================
$ws.Range("A26").Value = "attendances"
$ws.Range("A27").Value = "id"

$ws.Range("A28").Value = "class_id"
$ws.Range("A29").Value = "period"

# rows 28-29 use a new "middle row, green fill" style
$ws.Range("A28").Interior.Color = 5220458
$ws.Range("A28").Interior.PatternColor = 5220458
$ws.Range("A28").Borders.Item(7).LineStyle = 1
$ws.Range("A28").Borders.Item(10).LineStyle = 1
$ws.Range("A28").Copy()
$ws.Range("A29").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A29").Value = "period"

$ws.Range("A30").Value = "date"
# row 30 uses a new "last row, light/theme0 fill" style
$ws.Range("A9").Copy()
$ws.Range("A30").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A30").Value = "date"
$ws.Range("A30").Interior.ThemeColor = 2
$ws.Range("A30").Interior.PatternThemeColor = 2

# ===========================================================================
# attendance_class table (column A) - shifts up to rows 32-35
# ===========================================================================
$ws.Range("A31").Clear()

$ws.Range("A32").Value = "attendance_class"
$ws.Range("A16").Copy()
$ws.Range("A32").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A32").Value = "attendance_class"

$ws.Range("A33").Value = "attendance_id"
$ws.Range("A17").Copy()
$ws.Range("A33").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A33").Value = "attendance_id"

$ws.Range("A34").Value = "user_id (role học sinh)"
$ws.Range("A17").Copy()
$ws.Range("A34").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A34").Value = "user_id (role học sinh)"

$ws.Range("A35").Value = "is_present"

Write-Host "done"
